$d = $word.ActiveDocument

# --- Authors paragraph: add the missing contributor (Carol Bennett) ---
# Insert the previously-existing contributor text again right after "Authors:"
# (mirrors how the author re-typed/duplicated it ahead of the space run),
# and update the trailing contributor text to the new one.

$oldAuthors = "Warsame Yusuf, Methodologist, Ottawa Hospital Research Institute. Douglas G. Manuel, Senior Scientist & Physician, Ottawa Hospital Research Institute. Rostyslav Vyuha, Research Assistant, Ottawa Hospital Research Institute."
$newAuthors = "Warsame Yusuf, Methodologist, Ottawa Hospital Research Institute. Douglas G. Manuel, Senior Scientist & Physician, Ottawa Hospital Research Institute. Rostyslav Vyuha, Research Assistant, Ottawa Hospital Research Institute. Carol Bennett, Epidemiologist, Ottawa Hospital Research Institute."

$d.Content.Find.Execute($oldAuthors, $true, $false, $false, $false, $false, $true, 1, $false, $newAuthors, 2)

# --- Setting paragraph: no wording change, just touch it so Word recombines
#     the runs that were previously split mid-word ---
$d.Content.Find.Execute("The Canadian Community ", $true, $false, $false, $false, $false, $true, 1, $false, "The Canadian Community ", 2)

# --- Intervention paragraph: same, recombine the runs split mid-word ---
$d.Content.Find.Execute("We sought to use a cu", $true, $false, $false, $false, $false, $true, 1, $false, "We sought to use a cu", 2)

# --- Outcomes paragraph: recombine the runs split mid-word ---
$d.Content.Find.Execute("their names in p", $true, $false, $false, $false, $false, $true, 1, $false, "their names in p", 2)

# --- Implication paragraph: recombine the runs split mid-word ---
$d.Content.Find.Execute("generating a dataset of over ", $true, $false, $false, $false, $false, $true, 1, $false, "generating a dataset of over ", 2)
